# Fixed query issue for C3DC phs002599
# - Treatment query (row 5, "TreatmentTab"): drop the redundant CONCAT()
#   wrapper around REPLACE() in the "Treatment Agent" column expression.
# - Update the sheet's scroll/selection state to match where the author
#   left off editing (row 5 / cell C5) instead of row 7 / cell B7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$treatmentCell = $ws.Range("B5")
$oldQuery = $treatmentCell.Value()
$newQuery = $oldQuery.Replace(
    "CONCAT(REPLACE(trt.treatment_agent, ';', ', '))",
    "REPLACE(trt.treatment_agent, ';', ', ')"
)
$treatmentCell.Value = $newQuery

# Move the viewport / active selection to C5 (was B7), matching the
# author's final cursor position when they saved the workbook.
$ws.Range("C5").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
